# Update the marksheet "quiz" sheet with corrected Right/Total marks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Marking" row - right answer score per question (3 -> 5)
$ws.Range("B11").Value = 5

# "Total" row - total marks obtained (63 -> 105)
$ws.Range("B12").Value = 105

# "Total" row - Max column shows "correct/total" marks text (58/84 -> 105/140)
$ws.Range("E12").Value = "105/140"
